$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing row 48 (and below) down to 49
$ws.Rows.Item(48).Insert()

# Populate the new row 48 with the new data record
$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(48, 3).Value = "La Araucanía"
$ws.Cells.Item(48, 4).Value = 44799
$ws.Cells.Item(48, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(48, 5).Value = 9
$ws.Cells.Item(48, 6).Value = 100112035
$ws.Cells.Item(48, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 50
$ws.Cells.Item(48, 11).Value = 24000
$ws.Cells.Item(48, 12).Value = 25000
$ws.Cells.Item(48, 13).Value = 24600
$ws.Cells.Item(48, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(48, 15).Value = "Región Metropolitana"
$ws.Cells.Item(48, 16).Value = 2460
$ws.Cells.Item(48, 17).Value = 10
$ws.Cells.Item(48, 18).Value = "Hortaliza"
